# Updates crypto price/volume figures (and swaps the MXToken/WEMIXToken rows)
# per the Sun Oct 15 06:41:56 UTC 2023 GitHub Actions data refresh.
#
# Many "Price" column values are plain decimals (e.g. "1.00", "22.04") that
# Excel would otherwise auto-coerce to numbers on a plain .Value assignment,
# losing the original text formatting/representation. Forcing NumberFormat to
# Text ("@") before the assignment keeps them as strings; resetting the style
# back to "Normal" afterwards avoids leaving a stray text-format style behind
# so only the cell *values* change, matching the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "26.985.83"
Set-TextValue "E2" "  +0.22%  "
Set-TextValue "D3" "1.559.02"
Set-TextValue "E3" "  +0.55%  "
Set-TextValue "D5" "208.04"
Set-TextValue "E5" "  +0.71%  "
Set-TextValue "E6" "  +0.67%  "
Set-TextValue "E7" "  +0.18%  "
Set-TextValue "D8" "22.04"
Set-TextValue "E8" "  -0.04%  "
Set-TextValue "E9" "  +0.68%  "
Set-TextValue "E10" "  +1.87%  "
Set-TextValue "E11" "  -0.19%  "
Set-TextValue "D12" "1.780.59"
Set-TextValue "E12" "  +0.49%  "
Set-TextValue "D13" "1.562.85"
Set-TextValue "E13" "  +1.14%  "
Set-TextValue "D14" "3.72"
Set-TextValue "E14" "  -0.31%  "
Set-TextValue "E15" "  +0.11%  "
Set-TextValue "D16" "26.995.96"
Set-TextValue "E16" "  +0.26%  "
Set-TextValue "D17" "61.76"
Set-TextValue "E17" "  +0.24%  "
Set-TextValue "D18" "0.0₃0706"
Set-TextValue "E18" "  +1.47%  "
Set-TextValue "D19" "215.46"
Set-TextValue "E19" "  -0.89%  "
Set-TextValue "E20" "  +1.49%  "
Set-TextValue "E22" "  +2.20%  "
Set-TextValue "E23" "  +0.16%  "
Set-TextValue "D24" "1.94"
Set-TextValue "E24" "  -1.01%  "
Set-TextValue "D25" "152.73"
Set-TextValue "E25" "  -0.90%  "
Set-TextValue "E26" "  -0.21%  "
Set-TextValue "E28" "  +1.53%  "
Set-TextValue "D29" "1.01"
Set-TextValue "E29" "  +0.34%  "
Set-TextValue "D30" "0.0474"
Set-TextValue "E30" "  +1.40%  "
Set-TextValue "E31" "  +3.45%  "
Set-TextValue "E32" "  +0.35%  "
Set-TextValue "D33" "3.17"
Set-TextValue "E33" "  +3.37%  "
Set-TextValue "D34" "1.422.44"
Set-TextValue "E34" "  +0.02%  "
Set-TextValue "E35" "  +10.56%  "
Set-TextValue "E36" "  +1.14%  "
Set-TextValue "D37" "2.34"
Set-TextValue "E37" "  +2.31%  "
Set-TextValue "E38" "  +1.03%  "
Set-TextValue "E39" "  +1.99%  "
Set-TextValue "D40" "5.79"
Set-TextValue "E40" "  +0.63%  "
Set-TextValue "E41" "  +0.26%  "
Set-TextValue "E42" "  +0.19%  "
Set-TextValue "B43" "MXToken"
Set-TextValue "C43" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D43" "2.31"
Set-TextValue "E43" "  -0.09%  "
Set-TextValue "B44" "WEMIXToken"
Set-TextValue "C44" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D44" "0.999"
Set-TextValue "E44" "  +0.62%  "
Set-TextValue "D45" "64.59"
Set-TextValue "E45" "  +0.46%  "
Set-TextValue "E46" "  -1.18%  "
Set-TextValue "D47" "1.695.64"
Set-TextValue "E47" "  +0.55%  "
Set-TextValue "D48" "86.64"
Set-TextValue "E48" "  -1.19%  "
Set-TextValue "E49" "  +2.86%  "
Set-TextValue "D50" "0.0517"
Set-TextValue "E50" "  -0.27%  "
Set-TextValue "D51" "0.0959"
Set-TextValue "E51" "  +0.96%  "
